# Update the "想去人数" (interested-count) column F values on the
# "展览" and "全部类型" sheets, which mirror each other.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    2  = 1215
    3  = 69
    4  = 1126
    5  = 2464
    6  = 8017
    7  = 951
    8  = 493
    9  = 435
    11 = 461
    12 = 19
    13 = 183
    14 = 8387
    16 = 1461
    17 = 170
    20 = 209
    21 = 363
    22 = 211
    23 = 168
    28 = 1189
    29 = 87
    32 = 77
    34 = 51
    35 = 92
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
